# Apply updated distribution / cumulative / threshold-crossing values
# per the commit "Added Tire Type Filtering ... cleanup of Tire Type extraction"
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Step1_Data")
$ws1.Range("D2").Value = 0.00846908235109697
$ws1.Range("E2").Value = 0.3188526872842686
$ws1.Range("F2").Value = 0.2077865211829298
$ws1.Range("G2").Value = 0.08742545499143875
$ws1.Range("I2").Value = 0.09129936990231718
$ws1.Range("K2").Value = 0.02334904656840799
$ws1.Range("L2").Value = 0.007103925862645527
$ws1.Range("P2").Value = 0.1068098464846829
$ws1.Range("R2").Value = 0.120417380018447
$ws1.Range("T2").Value = 0.0110851649931781
$ws1.Range("U2").Value = 0.01193216943766957
$ws1.Range("V2").Value = 0.003301717500967987
$ws1.Range("Z2").Value = 0.002167633421949566
$ws1.Range("D3").Value = 0.3767972457842637
$ws1.Range("E3").Value = 0.05730741505477022
$ws1.Range("F3").Value = 0.1874661622207227
$ws1.Range("H3").Value = 0.05028953164038577
$ws1.Range("I3").Value = 0.01430633411554578
$ws1.Range("J3").Value = 0.02891949993715875
$ws1.Range("M3").Value = 0.02909824278344373
$ws1.Range("O3").Value = 0.07808529692760699
$ws1.Range("P3").Value = 0.05873596993292772
$ws1.Range("Q3").Value = 0.05597189309646743
$ws1.Range("R3").Value = 0.0008295373899893587
$ws1.Range("S3").Value = 0.04386980213841559
$ws1.Range("U3").Value = 0.01178174454897287
$ws1.Range("X3").Value = 0.002336714839298201
$ws1.Range("AD3").Value = 0.004204609590031098
$ws1.Range("D4").Value = 0.04988328848184409
$ws1.Range("E4").Value = 0.07000401020138369
$ws1.Range("F4").Value = 0.3360009134340134
$ws1.Range("G4").Value = 0.02445535884351198
$ws1.Range("H4").Value = 0.01402715727631138
$ws1.Range("I4").Value = 0.03234376884128307
$ws1.Range("K4").Value = 0.03421934393457213
$ws1.Range("M4").Value = 0.02365514701520149
$ws1.Range("N4").Value = 0.001037573127988074
$ws1.Range("O4").Value = 0.01412041767665762
$ws1.Range("P4").Value = 0.2047407026575661
$ws1.Range("R4").Value = 0.06830688201614321
$ws1.Range("S4").Value = 0.1022757003310224
$ws1.Range("U4").Value = 0.01400824204286441
$ws1.Range("AE4").Value = 0.007657359152643557
$ws1.Range("AF4").Value = 0.003264134966993256
$ws1.Range("E5").Value = 0.4380071043452757
$ws1.Range("F5").Value = 0.00645876020859655
$ws1.Range("G5").Value = 0.2114089636925653
$ws1.Range("I5").Value = 0.04120369649483067
$ws1.Range("J5").Value = 0.01787193087010551
$ws1.Range("K5").Value = 0.01388636553372598
$ws1.Range("N5").Value = 0.02817252535933436
$ws1.Range("P5").Value = 0.02358377239963174
$ws1.Range("Q5").Value = 0.07605672664937278
$ws1.Range("R5").Value = 0.02486100362702198
$ws1.Range("S5").Value = 0.04444981017371454
$ws1.Range("T5").Value = 0.01206184682424093
$ws1.Range("V5").Value = 0.03225606474143606
$ws1.Range("Y5").Value = 0.02795287080593106
$ws1.Range("AB5").Value = 0.0007109193021792769
$ws1.Range("AE5").Value = 0.001057638972037749
$ws1.Range("D6").Value = 0.23857682384639
$ws1.Range("E6").Value = 0.1183177961080697
$ws1.Range("F6").Value = 0.08334951731064449
$ws1.Range("H6").Value = 0.09449887568067798
$ws1.Range("K6").Value = 0.03542897972147265
$ws1.Range("M6").Value = 0.0197009116047755
$ws1.Range("N6").Value = 0.04573362731079309
$ws1.Range("O6").Value = 0.02666005970610753
$ws1.Range("P6").Value = 0.04327400050540442
$ws1.Range("Q6").Value = 0.07869728635107162
$ws1.Range("S6").Value = 0.1696427746041301
$ws1.Range("T6").Value = 0.01757004225717736
$ws1.Range("U6").Value = 0.02818718048113442
$ws1.Range("AF6").Value = 0.0003621245121511873

$ws2 = $wb.Worksheets.Item("Step2_Sj")
$ws2.Range("D2").Value = 0.00846908235109697
$ws2.Range("E2").Value = 0.3273217696353655
$ws2.Range("F2").Value = 0.5351082908182954
$ws2.Range("G2").Value = 0.6225337458097341
$ws2.Range("H2").Value = 0.6225337458097341
$ws2.Range("I2").Value = 0.7138331157120512
$ws2.Range("J2").Value = 0.7138331157120512
$ws2.Range("K2").Value = 0.7371821622804592
$ws2.Range("L2").Value = 0.7442860881431047
$ws2.Range("M2").Value = 0.7442860881431047
$ws2.Range("N2").Value = 0.7442860881431047
$ws2.Range("O2").Value = 0.7442860881431047
$ws2.Range("P2").Value = 0.8510959346277877
$ws2.Range("Q2").Value = 0.8510959346277877
$ws2.Range("R2").Value = 0.9715133146462347
$ws2.Range("S2").Value = 0.9715133146462347
$ws2.Range("T2").Value = 0.9825984796394128
$ws2.Range("U2").Value = 0.9945306490770824
$ws2.Range("V2").Value = 0.9978323665780504
$ws2.Range("W2").Value = 0.9978323665780504
$ws2.Range("X2").Value = 0.9978323665780504
$ws2.Range("Y2").Value = 0.9978323665780504
$ws2.Range("Z2").Value = 1
$ws2.Range("AA2").Value = 1
$ws2.Range("AB2").Value = 1
$ws2.Range("AC2").Value = 1
$ws2.Range("AD2").Value = 1
$ws2.Range("AE2").Value = 1
$ws2.Range("AF2").Value = 1
$ws2.Range("AG2").Value = 1
$ws2.Range("AH2").Value = 1
$ws2.Range("AI2").Value = 1
$ws2.Range("AJ2").Value = 1
$ws2.Range("D3").Value = 0.3767972457842637
$ws2.Range("E3").Value = 0.4341046608390339
$ws2.Range("F3").Value = 0.6215708230597566
$ws2.Range("G3").Value = 0.6215708230597566
$ws2.Range("H3").Value = 0.6718603547001424
$ws2.Range("I3").Value = 0.6861666888156882
$ws2.Range("J3").Value = 0.715086188752847
$ws2.Range("K3").Value = 0.715086188752847
$ws2.Range("L3").Value = 0.715086188752847
$ws2.Range("M3").Value = 0.7441844315362908
$ws2.Range("N3").Value = 0.7441844315362908
$ws2.Range("O3").Value = 0.8222697284638978
$ws2.Range("P3").Value = 0.8810056983968255
$ws2.Range("Q3").Value = 0.936977591493293
$ws2.Range("R3").Value = 0.9378071288832823
$ws2.Range("S3").Value = 0.981676931021698
$ws2.Range("T3").Value = 0.981676931021698
$ws2.Range("U3").Value = 0.9934586755706708
$ws2.Range("V3").Value = 0.9934586755706708
$ws2.Range("W3").Value = 0.9934586755706708
$ws2.Range("X3").Value = 0.995795390409969
$ws2.Range("Y3").Value = 0.995795390409969
$ws2.Range("Z3").Value = 0.995795390409969
$ws2.Range("AA3").Value = 0.995795390409969
$ws2.Range("AB3").Value = 0.995795390409969
$ws2.Range("AC3").Value = 0.995795390409969
$ws2.Range("D4").Value = 0.04988328848184409
$ws2.Range("E4").Value = 0.1198872986832278
$ws2.Range("F4").Value = 0.4558882121172412
$ws2.Range("G4").Value = 0.4803435709607532
$ws2.Range("H4").Value = 0.4943707282370646
$ws2.Range("I4").Value = 0.5267144970783476
$ws2.Range("J4").Value = 0.5267144970783476
$ws2.Range("K4").Value = 0.5609338410129198
$ws2.Range("L4").Value = 0.5609338410129198
$ws2.Range("M4").Value = 0.5845889880281212
$ws2.Range("N4").Value = 0.5856265611561092
$ws2.Range("O4").Value = 0.5997469788327668
$ws2.Range("P4").Value = 0.804487681490333
$ws2.Range("Q4").Value = 0.804487681490333
$ws2.Range("R4").Value = 0.8727945635064762
$ws2.Range("S4").Value = 0.9750702638374986
$ws2.Range("T4").Value = 0.9750702638374986
$ws2.Range("U4").Value = 0.9890785058803631
$ws2.Range("V4").Value = 0.9890785058803631
$ws2.Range("W4").Value = 0.9890785058803631
$ws2.Range("X4").Value = 0.9890785058803631
$ws2.Range("Y4").Value = 0.9890785058803631
$ws2.Range("Z4").Value = 0.9890785058803631
$ws2.Range("AA4").Value = 0.9890785058803631
$ws2.Range("AB4").Value = 0.9890785058803631
$ws2.Range("AC4").Value = 0.9890785058803631
$ws2.Range("AD4").Value = 0.9890785058803631
$ws2.Range("AE4").Value = 0.9967358650330066
$ws2.Range("AF4").Value = 0.9999999999999999
$ws2.Range("AG4").Value = 0.9999999999999999
$ws2.Range("AH4").Value = 0.9999999999999999
$ws2.Range("AI4").Value = 0.9999999999999999
$ws2.Range("AJ4").Value = 0.9999999999999999
$ws2.Range("E5").Value = 0.4380071043452757
$ws2.Range("F5").Value = 0.4444658645538722
$ws2.Range("G5").Value = 0.6558748282464375
$ws2.Range("H5").Value = 0.6558748282464375
$ws2.Range("I5").Value = 0.6970785247412682
$ws2.Range("J5").Value = 0.7149504556113737
$ws2.Range("K5").Value = 0.7288368211450997
$ws2.Range("L5").Value = 0.7288368211450997
$ws2.Range("M5").Value = 0.7288368211450997
$ws2.Range("N5").Value = 0.7570093465044341
$ws2.Range("O5").Value = 0.7570093465044341
$ws2.Range("P5").Value = 0.7805931189040658
$ws2.Range("Q5").Value = 0.8566498455534386
$ws2.Range("R5").Value = 0.8815108491804606
$ws2.Range("S5").Value = 0.9259606593541752
$ws2.Range("T5").Value = 0.9380225061784161
$ws2.Range("U5").Value = 0.9380225061784161
$ws2.Range("V5").Value = 0.9702785709198521
$ws2.Range("W5").Value = 0.9702785709198521
$ws2.Range("X5").Value = 0.9702785709198521
$ws2.Range("Y5").Value = 0.9982314417257832
$ws2.Range("Z5").Value = 0.9982314417257832
$ws2.Range("AA5").Value = 0.9982314417257832
$ws2.Range("AB5").Value = 0.9989423610279624
$ws2.Range("AC5").Value = 0.9989423610279624
$ws2.Range("AD5").Value = 0.9989423610279624
$ws2.Range("D6").Value = 0.23857682384639
$ws2.Range("E6").Value = 0.3568946199544597
$ws2.Range("F6").Value = 0.4402441372651042
$ws2.Range("G6").Value = 0.4402441372651042
$ws2.Range("H6").Value = 0.5347430129457822
$ws2.Range("I6").Value = 0.5347430129457822
$ws2.Range("J6").Value = 0.5347430129457822
$ws2.Range("K6").Value = 0.5701719926672548
$ws2.Range("L6").Value = 0.5701719926672548
$ws2.Range("M6").Value = 0.5898729042720303
$ws2.Range("N6").Value = 0.6356065315828234
$ws2.Range("O6").Value = 0.662266591288931
$ws2.Range("P6").Value = 0.7055405917943354
$ws2.Range("Q6").Value = 0.784237878145407
$ws2.Range("R6").Value = 0.784237878145407
$ws2.Range("S6").Value = 0.9538806527495372
$ws2.Range("T6").Value = 0.9714506950067145
$ws2.Range("U6").Value = 0.999637875487849
$ws2.Range("V6").Value = 0.999637875487849
$ws2.Range("W6").Value = 0.999637875487849
$ws2.Range("X6").Value = 0.999637875487849
$ws2.Range("Y6").Value = 0.999637875487849
$ws2.Range("Z6").Value = 0.999637875487849
$ws2.Range("AA6").Value = 0.999637875487849
$ws2.Range("AB6").Value = 0.999637875487849
$ws2.Range("AC6").Value = 0.999637875487849
$ws2.Range("AD6").Value = 0.999637875487849
$ws2.Range("AE6").Value = 0.999637875487849
$ws2.Range("AF6").Value = 1
$ws2.Range("AG6").Value = 1
$ws2.Range("AH6").Value = 1
$ws2.Range("AI6").Value = 1
$ws2.Range("AJ6").Value = 1

$ws3 = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws3.Range("F2").Value = 0.5351082908182954
$ws3.Range("D3").Value = 5
$ws3.Range("F3").Value = 0.6215708230597566
$ws3.Range("G3").Value = 4
$ws3.Range("D4").Value = 8
$ws3.Range("F4").Value = 0.5267144970783476
$ws3.Range("G4").Value = 6
$ws3.Range("D5").Value = 6
$ws3.Range("F5").Value = 0.6558748282464375
$ws3.Range("G5").Value = 4
$ws3.Range("D6").Value = 7
$ws3.Range("F6").Value = 0.5347430129457822
$ws3.Range("G6").Value = 6

$ws4 = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws4.Range("F2").Value = 0.7138331157120512
$ws4.Range("D3").Value = 9
$ws4.Range("F3").Value = 0.715086188752847
$ws4.Range("G3").Value = 8
$ws4.Range("F4").Value = 0.804487681490333
$ws4.Range("D5").Value = 9
$ws4.Range("F5").Value = 0.7149504556113737
$ws4.Range("G5").Value = 7
$ws4.Range("F6").Value = 0.7055405917943354

$ws5 = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws5.Range("F2").Value = 0.8510959346277877
$ws5.Range("D3").Value = 14
$ws5.Range("F3").Value = 0.8222697284638978
$ws5.Range("G3").Value = 13
$ws5.Range("F4").Value = 0.804487681490333
$ws5.Range("D5").Value = 16
$ws5.Range("F5").Value = 0.8566498455534386
$ws5.Range("G5").Value = 14
$ws5.Range("F6").Value = 0.9538806527495372

$ws6 = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws6.Range("F2").Value = 0.9715133146462347
$ws6.Range("D3").Value = 16
$ws6.Range("F3").Value = 0.936977591493293
$ws6.Range("G3").Value = 15
$ws6.Range("F4").Value = 0.9750702638374986
$ws6.Range("D5").Value = 18
$ws6.Range("F5").Value = 0.9259606593541752
$ws6.Range("G5").Value = 16
$ws6.Range("F6").Value = 0.9538806527495372

